$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

Set-TextValue $ws.Range('D2') '28.034.79'
Set-TextValue $ws.Range('E2') '  +2.14%  '
Set-TextValue $ws.Range('D3') '1.643.65'
Set-TextValue $ws.Range('E3') '  +0.40%  '
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.29%  '
Set-TextValue $ws.Range('D5') '212.72'
Set-TextValue $ws.Range('E5') '  +0.16%  '
Set-TextValue $ws.Range('E6') '  -1.39%  '
Set-TextValue $ws.Range('D7') '0.998'
Set-TextValue $ws.Range('E7') '  -0.44%  '
Set-TextValue $ws.Range('D8') '23.43'
Set-TextValue $ws.Range('E8') '  +1.30%  '
Set-TextValue $ws.Range('E9') '  +3.14%  '
Set-TextValue $ws.Range('D10') '0.0613'
Set-TextValue $ws.Range('E10') '  +0.62%  '
Set-TextValue $ws.Range('E11') '  +0.55%  '
Set-TextValue $ws.Range('D12') '1.874.29'
Set-TextValue $ws.Range('E12') '  +0.20%  '
Set-TextValue $ws.Range('D13') '1.627.91'
Set-TextValue $ws.Range('E13') '  -0.69%  '
Set-TextValue $ws.Range('D14') '4.06'
Set-TextValue $ws.Range('E14') '  +1.40%  '
Set-TextValue $ws.Range('D15') '0.562'
Set-TextValue $ws.Range('E15') '  -3.26%  '
Set-TextValue $ws.Range('D16') '64.73'
Set-TextValue $ws.Range('E16') '  +0.98%  '
Set-TextValue $ws.Range('D17') '27.987.62'
Set-TextValue $ws.Range('E17') '  +2.00%  '
Set-TextValue $ws.Range('D18') '233.89'
Set-TextValue $ws.Range('E18') '  +2.24%  '
Set-TextValue $ws.Range('D19') '0.0₃0724'
Set-TextValue $ws.Range('E19') '  +0.33%  '
Set-TextValue $ws.Range('D20') '7.65'
Set-TextValue $ws.Range('E20') '  +1.71%  '
Set-TextValue $ws.Range('E21') '  -0.25%  '
Set-TextValue $ws.Range('D22') '4.32'
Set-TextValue $ws.Range('E22') '  +0.68%  '
Set-TextValue $ws.Range('D23') '10.04'
Set-TextValue $ws.Range('E23') '  +3.95%  '
Set-TextValue $ws.Range('E24') '  +4.76%  '
Set-TextValue $ws.Range('D25') '150.58'
Set-TextValue $ws.Range('E25') '  +0.96%  '
Set-TextValue $ws.Range('D26') '6.95'
Set-TextValue $ws.Range('E26') '  -0.40%  '
Set-TextValue $ws.Range('E27') '  -0.56%  '
Set-TextValue $ws.Range('E28') '  +1.23%  '
Set-TextValue $ws.Range('D29') '0.999'
Set-TextValue $ws.Range('E29') '  -0.33%  '
Set-TextValue $ws.Range('E30') '  +0.32%  '
Set-TextValue $ws.Range('E31') '  -0.94%  '
Set-TextValue $ws.Range('D32') '3.32'
Set-TextValue $ws.Range('E32') '  +0.97%  '
Set-TextValue $ws.Range('D33') '1.474.11'
Set-TextValue $ws.Range('E33') '  +4.18%  '
Set-TextValue $ws.Range('E34') '  -1.86%  '
Set-TextValue $ws.Range('D35') '1.56'
Set-TextValue $ws.Range('E35') '  -2.35%  '
Set-TextValue $ws.Range('E36') '  -0.31%  '
Set-TextValue $ws.Range('E37') '  -0.34%  '
Set-TextValue $ws.Range('D38') '0.883'
Set-TextValue $ws.Range('E38') '  +0.65%  '
Set-TextValue $ws.Range('E39') '  +0.92%  '
Set-TextValue $ws.Range('D40') '0.919'
Set-TextValue $ws.Range('E40') '  +13.13%  '
Set-TextValue $ws.Range('D41') '69.82'
Set-TextValue $ws.Range('E41') '  +7.91%  '
Set-TextValue $ws.Range('E42') '  -0.25%  '
Set-TextValue $ws.Range('E43') '  -1.68%  '
Set-TextValue $ws.Range('D44') '2.45'
Set-TextValue $ws.Range('E44') '  -2.34%  '
Set-TextValue $ws.Range('E45') '  +0.07%  '
Set-TextValue $ws.Range('E46') '  -1.14%  '
Set-TextValue $ws.Range('D47') '1.785.39'
Set-TextValue $ws.Range('E47') '  +0.32%  '
Set-TextValue $ws.Range('E48') '  +2.77%  '
Set-TextValue $ws.Range('D49') '86.85'
Set-TextValue $ws.Range('E49') '  +1.43%  '
Set-TextValue $ws.Range('E50') '  +0.22%  '
Set-TextValue $ws.Range('D51') '0.0994'
Set-TextValue $ws.Range('E51') '  +0.25%  '
